# Apply the "changes to excel sheet" commit:
#  - Register sheet: remove the sample data rows (2,3,5,6) and clear the
#    remaining placeholder row (row 4), leaving only the header row and an
#    empty, still-styled row 4.
#  - SignIn sheet: clear the sample data row (row 2), leaving only the
#    header row and an empty, still-styled row 2.
#  - PythonCode sheet content is untouched (its shared-string indices shift
#    automatically once the now-unused strings from Register/SignIn are
#    dropped by the engine on save).
#  - Update the active sheet/selection so SignIn becomes the active tab,
#    matching the saved view state in the target workbook.

$wb = $excel.ActiveWorkbook

$register = $wb.Worksheets.Item("Register")
$signIn   = $wb.Worksheets.Item("SignIn")

# --- Register ("Register" sheet) -------------------------------------------------
# Remove the ninjatesters / Team104 / Numpy / kodetesters sample rows entirely.
$register.Range("A2:C3").ClearContents()
$register.Range("A5:C6").ClearContents()
# Keep row 4 (with its style) but blank out its values, including A4's text.
$register.Range("A4:C4").ClearContents()

# --- SignIn sheet ------------------------------------------------------------
# Remove the kodetesters / numpyninja24 sample row, keeping B2's style.
$signIn.Range("A2:B2").ClearContents()

# --- Selections / active sheet ------------------------------------------------
$register.Activate() | Out-Null
$register.Range("A2:C6").Select() | Out-Null

$signIn.Activate() | Out-Null
$signIn.Range("A2:B2").Select() | Out-Null
